$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Hunk 1: add "Minigames" bullet text to the existing empty bullet under
# "Transmission Grid Operator", then add three new sub-bullets
# ("Memory", "Hangman", and an empty one carrying the relocated
# "_GoBack" bookmark) right after it.
# ---------------------------------------------------------------------

$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Inter-grid connectivity") {
        $targetPara = $p.Next()
    }
}

$r = $targetPara.Range
$r.InsertAfter("Minigames")
$r.LanguageID = "en-US"
$targetPara.Range.InsertParagraphAfter()

$memoryPara = $targetPara.Next()
$memoryPara.Range.ListFormat.ListLevelNumber = 2
$memoryPara.Range.InsertBefore("Memory")
$memoryPara.Range.InsertParagraphAfter()

$hangmanPara = $memoryPara.Next()
$hangmanPara.Range.ListFormat.ListLevelNumber = 2
$hangmanPara.Range.InsertBefore("Hangman")
$hangmanPara.Range.InsertParagraphAfter()

$bmPara = $hangmanPara.Next()
$bmPara.Range.ListFormat.ListLevelNumber = 2

# Bookmarks.Add can't anchor into a wholly run-less (empty) paragraph,
# so stage a throwaway character, plant the bookmark, then remove the
# character again -- the bookmark (a zero-length point) survives.
$phStart = $bmPara.Range.Start
$ph = $d.Range($phStart, $phStart)
$ph.InsertAfter("X")

$bmRange = $d.Range($phStart, $phStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($phStart, $phStart + 1).Delete()

# ---------------------------------------------------------------------
# Hunk 2: drop the stale <w:lastRenderedPageBreak/> hint in front of the
# "Andere" heading -- re-stamping the run's text clears it.
# ---------------------------------------------------------------------

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Andere" -and $p.Style.NameLocal -eq "Heading 1") {
        $p.Range.Text = "Andere"
    }
}

Write-Host "Done"
